# Add a new trade record (row 16) to the sheet, mirroring the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 9130.19
$ws.Range("B16").Value = 9402.8700000000008
$ws.Range("C16").Value = 294.14
$ws.Range("D16").Value = 302.66000000000003
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = 2.9
$ws.Range("G16").Value = 42626.544340277775
$ws.Range("G16").NumberFormat = "m/d/yy h:mm"
$ws.Range("H16").Value = $false
